$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 39-112: a new week of Perejil price data was
# inserted into the middle of the table on 5 separate dates, pushing the
# previously-existing rows down. Re-write each affected row in place with
# its new (shifted) content. ---

$ws.Range("D39").Value2 = 44565
$ws.Range("J39").Value2 = 3000
$ws.Range("K39").Value2 = 3000
$ws.Range("L39").Value2 = 3500
$ws.Range("M39").Value2 = 3250
$ws.Range("P39").Value2 = 2167
$ws.Range("D40").Value2 = 44236
$ws.Range("J40").Value2 = 3200
$ws.Range("K40").Value2 = 1500
$ws.Range("L40").Value2 = 2000
$ws.Range("M40").Value2 = 1750
$ws.Range("P40").Value2 = 1167
$ws.Range("D41").Value2 = 44546
$ws.Range("J41").Value2 = 2800
$ws.Range("K41").Value2 = 2000
$ws.Range("L41").Value2 = 2500
$ws.Range("M41").Value2 = 2250
$ws.Range("P41").Value2 = 1500
$ws.Range("D42").Value2 = 44222
$ws.Range("J42").Value2 = 2800
$ws.Range("K42").Value2 = 1300
$ws.Range("L42").Value2 = 1500
$ws.Range("M42").Value2 = 1400
$ws.Range("P42").Value2 = 933
$ws.Range("D43").Value2 = 44488
$ws.Range("J43").Value2 = 3000
$ws.Range("K43").Value2 = 1300
$ws.Range("L43").Value2 = 1500
$ws.Range("M43").Value2 = 1400
$ws.Range("P43").Value2 = 933
$ws.Range("D44").Value2 = 44446
$ws.Range("J44").Value2 = 3200
$ws.Range("K44").Value2 = 2000
$ws.Range("L44").Value2 = 2500
$ws.Range("M44").Value2 = 2250
$ws.Range("P44").Value2 = 1500
$ws.Range("D45").Value2 = 44483
$ws.Range("J45").Value2 = 3300
$ws.Range("K45").Value2 = 1500
$ws.Range("L45").Value2 = 2000
$ws.Range("M45").Value2 = 1750
$ws.Range("P45").Value2 = 1167
$ws.Range("D46").Value2 = 44215
$ws.Range("J46").Value2 = 2800
$ws.Range("K46").Value2 = 1300
$ws.Range("L46").Value2 = 1500
$ws.Range("M46").Value2 = 1400
$ws.Range("P46").Value2 = 933
$ws.Range("D47").Value2 = 44411
$ws.Range("J47").Value2 = 3300
$ws.Range("K47").Value2 = 2000
$ws.Range("L47").Value2 = 2500
$ws.Range("M47").Value2 = 2250
$ws.Range("P47").Value2 = 1500
$ws.Range("D48").Value2 = 44469
$ws.Range("J48").Value2 = 3100
$ws.Range("K48").Value2 = 1500
$ws.Range("L48").Value2 = 2000
$ws.Range("M48").Value2 = 1750
$ws.Range("P48").Value2 = 1167
$ws.Range("D49").Value2 = 44509
$ws.Range("J49").Value2 = 2940
$ws.Range("K49").Value2 = 1300
$ws.Range("L49").Value2 = 1500
$ws.Range("M49").Value2 = 1400
$ws.Range("P49").Value2 = 933
$ws.Range("D50").Value2 = 44348
$ws.Range("J50").Value2 = 3360
$ws.Range("K50").Value2 = 1300
$ws.Range("L50").Value2 = 1500
$ws.Range("M50").Value2 = 1400
$ws.Range("P50").Value2 = 933
$ws.Range("D51").Value2 = 44224
$ws.Range("J51").Value2 = 2800
$ws.Range("K51").Value2 = 1300
$ws.Range("L51").Value2 = 1500
$ws.Range("M51").Value2 = 1400
$ws.Range("P51").Value2 = 933
$ws.Range("D52").Value2 = 44390
$ws.Range("J52").Value2 = 3200
$ws.Range("K52").Value2 = 1500
$ws.Range("L52").Value2 = 2000
$ws.Range("M52").Value2 = 1750
$ws.Range("P52").Value2 = 1167
$ws.Range("D53").Value2 = 44168
$ws.Range("J53").Value2 = 2800
$ws.Range("K53").Value2 = 1300
$ws.Range("L53").Value2 = 1500
$ws.Range("M53").Value2 = 1400
$ws.Range("P53").Value2 = 933
$ws.Range("D54").Value2 = 44397
$ws.Range("J54").Value2 = 3200
$ws.Range("K54").Value2 = 1500
$ws.Range("L54").Value2 = 2000
$ws.Range("M54").Value2 = 1750
$ws.Range("P54").Value2 = 1167
$ws.Range("D55").Value2 = 44308
$ws.Range("J55").Value2 = 3200
$ws.Range("K55").Value2 = 1300
$ws.Range("L55").Value2 = 1500
$ws.Range("M55").Value2 = 1400
$ws.Range("P55").Value2 = 933
$ws.Range("D56").Value2 = 44579
$ws.Range("J56").Value2 = 3000
$ws.Range("K56").Value2 = 2500
$ws.Range("L56").Value2 = 3000
$ws.Range("M56").Value2 = 2750
$ws.Range("P56").Value2 = 1833
$ws.Range("D57").Value2 = 44574
$ws.Range("J57").Value2 = 3200
$ws.Range("K57").Value2 = 3000
$ws.Range("L57").Value2 = 3500
$ws.Range("M57").Value2 = 3250
$ws.Range("P57").Value2 = 2167
$ws.Range("D58").Value2 = 44231
$ws.Range("J58").Value2 = 2800
$ws.Range("K58").Value2 = 1300
$ws.Range("L58").Value2 = 1500
$ws.Range("M58").Value2 = 1400
$ws.Range("P58").Value2 = 933
$ws.Range("D59").Value2 = 44435
$ws.Range("J59").Value2 = 6560
$ws.Range("K59").Value2 = 2000
$ws.Range("L59").Value2 = 2500
$ws.Range("M59").Value2 = 2250
$ws.Range("P59").Value2 = 1500
$ws.Range("D60").Value2 = 44336
$ws.Range("J60").Value2 = 3360
$ws.Range("K60").Value2 = 1300
$ws.Range("L60").Value2 = 1500
$ws.Range("M60").Value2 = 1400
$ws.Range("P60").Value2 = 933
$ws.Range("D61").Value2 = 44350
$ws.Range("J61").Value2 = 3340
$ws.Range("K61").Value2 = 1200
$ws.Range("L61").Value2 = 1500
$ws.Range("M61").Value2 = 1350
$ws.Range("P61").Value2 = 900
$ws.Range("D62").Value2 = 44413
$ws.Range("J62").Value2 = 3360
$ws.Range("K62").Value2 = 2000
$ws.Range("L62").Value2 = 2500
$ws.Range("M62").Value2 = 2250
$ws.Range("P62").Value2 = 1500
$ws.Range("D63").Value2 = 44329
$ws.Range("J63").Value2 = 3300
$ws.Range("K63").Value2 = 1300
$ws.Range("L63").Value2 = 1500
$ws.Range("M63").Value2 = 1400
$ws.Range("P63").Value2 = 933
$ws.Range("D64").Value2 = 44553
$ws.Range("J64").Value2 = 3200
$ws.Range("K64").Value2 = 2000
$ws.Range("L64").Value2 = 2500
$ws.Range("M64").Value2 = 2250
$ws.Range("P64").Value2 = 1500
$ws.Range("D65").Value2 = 44434
$ws.Range("J65").Value2 = 3360
$ws.Range("K65").Value2 = 2000
$ws.Range("L65").Value2 = 2500
$ws.Range("M65").Value2 = 2250
$ws.Range("P65").Value2 = 1500
$ws.Range("D66").Value2 = 44551
$ws.Range("J66").Value2 = 2880
$ws.Range("K66").Value2 = 2000
$ws.Range("L66").Value2 = 2500
$ws.Range("M66").Value2 = 2250
$ws.Range("P66").Value2 = 1500
$ws.Range("D67").Value2 = 44476
$ws.Range("J67").Value2 = 3000
$ws.Range("K67").Value2 = 1500
$ws.Range("L67").Value2 = 2000
$ws.Range("M67").Value2 = 1750
$ws.Range("P67").Value2 = 1167
$ws.Range("D68").Value2 = 44586
$ws.Range("J68").Value2 = 2800
$ws.Range("K68").Value2 = 2500
$ws.Range("L68").Value2 = 3000
$ws.Range("M68").Value2 = 2750
$ws.Range("P68").Value2 = 1833
$ws.Range("D69").Value2 = 44558
$ws.Range("J69").Value2 = 2900
$ws.Range("K69").Value2 = 2300
$ws.Range("L69").Value2 = 2500
$ws.Range("M69").Value2 = 2400
$ws.Range("P69").Value2 = 1600
$ws.Range("D70").Value2 = 44252
$ws.Range("J70").Value2 = 3600
$ws.Range("K70").Value2 = 1500
$ws.Range("L70").Value2 = 2000
$ws.Range("M70").Value2 = 1750
$ws.Range("P70").Value2 = 1167
$ws.Range("D71").Value2 = 44238
$ws.Range("J71").Value2 = 3200
$ws.Range("K71").Value2 = 1500
$ws.Range("L71").Value2 = 2000
$ws.Range("M71").Value2 = 1750
$ws.Range("P71").Value2 = 1167
$ws.Range("D72").Value2 = 44530
$ws.Range("J72").Value2 = 2800
$ws.Range("K72").Value2 = 1500
$ws.Range("L72").Value2 = 2000
$ws.Range("M72").Value2 = 1750
$ws.Range("P72").Value2 = 1167
$ws.Range("D73").Value2 = 44327
$ws.Range("J73").Value2 = 3400
$ws.Range("K73").Value2 = 1300
$ws.Range("L73").Value2 = 1500
$ws.Range("M73").Value2 = 1400
$ws.Range("P73").Value2 = 933
$ws.Range("D74").Value2 = 44455
$ws.Range("J74").Value2 = 3200
$ws.Range("K74").Value2 = 2000
$ws.Range("L74").Value2 = 2500
$ws.Range("M74").Value2 = 2250
$ws.Range("P74").Value2 = 1500
$ws.Range("D75").Value2 = 44159
$ws.Range("J75").Value2 = 2900
$ws.Range("K75").Value2 = 1000
$ws.Range("L75").Value2 = 1500
$ws.Range("M75").Value2 = 1250
$ws.Range("P75").Value2 = 833
$ws.Range("D76").Value2 = 44362
$ws.Range("J76").Value2 = 3200
$ws.Range("K76").Value2 = 1500
$ws.Range("L76").Value2 = 2000
$ws.Range("M76").Value2 = 1750
$ws.Range("P76").Value2 = 1167
$ws.Range("D77").Value2 = 44462
$ws.Range("J77").Value2 = 3200
$ws.Range("K77").Value2 = 1500
$ws.Range("L77").Value2 = 2000
$ws.Range("M77").Value2 = 1750
$ws.Range("P77").Value2 = 1167
$ws.Range("D78").Value2 = 44364
$ws.Range("J78").Value2 = 3200
$ws.Range("K78").Value2 = 1500
$ws.Range("L78").Value2 = 2000
$ws.Range("M78").Value2 = 1750
$ws.Range("P78").Value2 = 1167
$ws.Range("D79").Value2 = 44376
$ws.Range("J79").Value2 = 3200
$ws.Range("K79").Value2 = 1500
$ws.Range("L79").Value2 = 2000
$ws.Range("M79").Value2 = 1750
$ws.Range("P79").Value2 = 1167
$ws.Range("D80").Value2 = 44313
$ws.Range("J80").Value2 = 3200
$ws.Range("K80").Value2 = 1300
$ws.Range("L80").Value2 = 1500
$ws.Range("M80").Value2 = 1400
$ws.Range("P80").Value2 = 933
$ws.Range("D81").Value2 = 44572
$ws.Range("J81").Value2 = 2900
$ws.Range("K81").Value2 = 3000
$ws.Range("L81").Value2 = 3500
$ws.Range("M81").Value2 = 3250
$ws.Range("P81").Value2 = 2167
$ws.Range("D82").Value2 = 44516
$ws.Range("J82").Value2 = 3100
$ws.Range("K82").Value2 = 1300
$ws.Range("L82").Value2 = 1500
$ws.Range("M82").Value2 = 1400
$ws.Range("P82").Value2 = 933
$ws.Range("D83").Value2 = 44257
$ws.Range("J83").Value2 = 3600
$ws.Range("K83").Value2 = 2000
$ws.Range("L83").Value2 = 2500
$ws.Range("M83").Value2 = 2250
$ws.Range("P83").Value2 = 1500
$ws.Range("D84").Value2 = 44322
$ws.Range("J84").Value2 = 3320
$ws.Range("K84").Value2 = 1300
$ws.Range("L84").Value2 = 1500
$ws.Range("M84").Value2 = 1400
$ws.Range("P84").Value2 = 933
$ws.Range("D85").Value2 = 44217
$ws.Range("J85").Value2 = 2800
$ws.Range("K85").Value2 = 1300
$ws.Range("L85").Value2 = 1500
$ws.Range("M85").Value2 = 1400
$ws.Range("P85").Value2 = 933
$ws.Range("D86").Value2 = 44540
$ws.Range("J86").Value2 = 3000
$ws.Range("K86").Value2 = 1500
$ws.Range("L86").Value2 = 2000
$ws.Range("M86").Value2 = 1750
$ws.Range("P86").Value2 = 1167
$ws.Range("D87").Value2 = 44511
$ws.Range("J87").Value2 = 3360
$ws.Range("K87").Value2 = 1300
$ws.Range("L87").Value2 = 1500
$ws.Range("M87").Value2 = 1400
$ws.Range("P87").Value2 = 933
$ws.Range("D88").Value2 = 44166
$ws.Range("J88").Value2 = 2800
$ws.Range("K88").Value2 = 1300
$ws.Range("L88").Value2 = 1500
$ws.Range("M88").Value2 = 1400
$ws.Range("P88").Value2 = 933
$ws.Range("D89").Value2 = 44175
$ws.Range("J89").Value2 = 3000
$ws.Range("K89").Value2 = 1300
$ws.Range("L89").Value2 = 1500
$ws.Range("M89").Value2 = 1400
$ws.Range("P89").Value2 = 933
$ws.Range("D90").Value2 = 44203
$ws.Range("J90").Value2 = 2800
$ws.Range("K90").Value2 = 1300
$ws.Range("L90").Value2 = 1500
$ws.Range("M90").Value2 = 1400
$ws.Range("P90").Value2 = 933
$ws.Range("D91").Value2 = 44161
$ws.Range("J91").Value2 = 3100
$ws.Range("K91").Value2 = 1300
$ws.Range("L91").Value2 = 1500
$ws.Range("M91").Value2 = 1400
$ws.Range("P91").Value2 = 933
$ws.Range("D92").Value2 = 44581
$ws.Range("J92").Value2 = 3100
$ws.Range("K92").Value2 = 2500
$ws.Range("L92").Value2 = 3000
$ws.Range("M92").Value2 = 2750
$ws.Range("P92").Value2 = 1833
$ws.Range("D93").Value2 = 44504
$ws.Range("J93").Value2 = 3200
$ws.Range("K93").Value2 = 1300
$ws.Range("L93").Value2 = 1500
$ws.Range("M93").Value2 = 1400
$ws.Range("P93").Value2 = 933
$ws.Range("D94").Value2 = 44567
$ws.Range("J94").Value2 = 3200
$ws.Range("K94").Value2 = 3000
$ws.Range("L94").Value2 = 3500
$ws.Range("M94").Value2 = 3250
$ws.Range("P94").Value2 = 2167
$ws.Range("D95").Value2 = 44280
$ws.Range("J95").Value2 = 3000
$ws.Range("K95").Value2 = 2000
$ws.Range("L95").Value2 = 2500
$ws.Range("M95").Value2 = 2250
$ws.Range("P95").Value2 = 1500
$ws.Range("D96").Value2 = 44532
$ws.Range("J96").Value2 = 3260
$ws.Range("K96").Value2 = 1800
$ws.Range("L96").Value2 = 2000
$ws.Range("M96").Value2 = 1900
$ws.Range("P96").Value2 = 1267
$ws.Range("D97").Value2 = 44334
$ws.Range("J97").Value2 = 3440
$ws.Range("K97").Value2 = 1300
$ws.Range("L97").Value2 = 1500
$ws.Range("M97").Value2 = 1400
$ws.Range("P97").Value2 = 933
$ws.Range("D98").Value2 = 44187
$ws.Range("J98").Value2 = 3100
$ws.Range("K98").Value2 = 1400
$ws.Range("L98").Value2 = 1500
$ws.Range("M98").Value2 = 1450
$ws.Range("P98").Value2 = 967
$ws.Range("D99").Value2 = 44266
$ws.Range("J99").Value2 = 3600
$ws.Range("K99").Value2 = 2000
$ws.Range("L99").Value2 = 2500
$ws.Range("M99").Value2 = 2250
$ws.Range("P99").Value2 = 1500
$ws.Range("D100").Value2 = 44371
$ws.Range("J100").Value2 = 3300
$ws.Range("K100").Value2 = 1500
$ws.Range("L100").Value2 = 2000
$ws.Range("M100").Value2 = 1750
$ws.Range("P100").Value2 = 1167
$ws.Range("D101").Value2 = 44259
$ws.Range("J101").Value2 = 3400
$ws.Range("K101").Value2 = 2000
$ws.Range("L101").Value2 = 2500
$ws.Range("M101").Value2 = 2250
$ws.Range("P101").Value2 = 1500
$ws.Range("D102").Value2 = 44250
$ws.Range("J102").Value2 = 3400
$ws.Range("K102").Value2 = 1500
$ws.Range("L102").Value2 = 2000
$ws.Range("M102").Value2 = 1750
$ws.Range("P102").Value2 = 1167
$ws.Range("D103").Value2 = 44285
$ws.Range("J103").Value2 = 3400
$ws.Range("K103").Value2 = 2000
$ws.Range("L103").Value2 = 2500
$ws.Range("M103").Value2 = 2250
$ws.Range("P103").Value2 = 1500
$ws.Range("D104").Value2 = 44264
$ws.Range("J104").Value2 = 3600
$ws.Range("K104").Value2 = 2000
$ws.Range("L104").Value2 = 2500
$ws.Range("M104").Value2 = 2250
$ws.Range("P104").Value2 = 1500
$ws.Range("D105").Value2 = 44523
$ws.Range("J105").Value2 = 2800
$ws.Range("K105").Value2 = 1500
$ws.Range("L105").Value2 = 2000
$ws.Range("M105").Value2 = 1750
$ws.Range("P105").Value2 = 1167
$ws.Range("D106").Value2 = 44399
$ws.Range("J106").Value2 = 3320
$ws.Range("K106").Value2 = 1500
$ws.Range("L106").Value2 = 2000
$ws.Range("M106").Value2 = 1750
$ws.Range("P106").Value2 = 1167
$ws.Range("D107").Value2 = 44441
$ws.Range("J107").Value2 = 3200
$ws.Range("K107").Value2 = 2000
$ws.Range("L107").Value2 = 2500
$ws.Range("M107").Value2 = 2250
$ws.Range("P107").Value2 = 1500
$ws.Range("D108").Value2 = 44315
$ws.Range("J108").Value2 = 3120
$ws.Range("K108").Value2 = 1300
$ws.Range("L108").Value2 = 1500
$ws.Range("M108").Value2 = 1400
$ws.Range("P108").Value2 = 933
$ws.Range("D109").Value2 = 44278
$ws.Range("J109").Value2 = 3400
$ws.Range("K109").Value2 = 2000
$ws.Range("L109").Value2 = 2500
$ws.Range("M109").Value2 = 2250
$ws.Range("P109").Value2 = 1500
$ws.Range("D110").Value2 = 44453
$ws.Range("J110").Value2 = 3200
$ws.Range("K110").Value2 = 2000
$ws.Range("L110").Value2 = 2500
$ws.Range("M110").Value2 = 2250
$ws.Range("P110").Value2 = 1500
$ws.Range("D111").Value2 = 44474
$ws.Range("J111").Value2 = 2800
$ws.Range("K111").Value2 = 1500
$ws.Range("L111").Value2 = 2000
$ws.Range("M111").Value2 = 1750
$ws.Range("P111").Value2 = 1167
$ws.Range("D112").Value2 = 44560
$ws.Range("J112").Value2 = 3400
$ws.Range("K112").Value2 = 2500
$ws.Range("L112").Value2 = 3000
$ws.Range("M112").Value2 = 2750
$ws.Range("P112").Value2 = 1833

# --- Append 5 new rows (113-117) at the bottom of the table: the rows that
# were pushed out of the 39-112 window by the inserted data above. ---

$ws.Range("A113").Value2 = 8
$ws.Range("B113").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C113").Value2 = "Coquimbo"
$ws.Range("E113").Value2 = 4
$ws.Range("F113").Value2 = 100112044
$ws.Range("G113").Value2 = "Perejil"
$ws.Range("H113").Value2 = "Sin especificar"
$ws.Range("I113").Value2 = "Primera"
$ws.Range("N113").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O113").Value2 = "Provincia del Elquí"
$ws.Range("Q113").Value2 = 1.5
$ws.Range("R113").Value2 = "Hortaliza"
$ws.Range("D113").Value2 = 44385
$ws.Range("D113").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J113").Value2 = 3320
$ws.Range("K113").Value2 = 1500
$ws.Range("L113").Value2 = 2000
$ws.Range("M113").Value2 = 1750
$ws.Range("P113").Value2 = 1167

$ws.Range("A114").Value2 = 8
$ws.Range("B114").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C114").Value2 = "Coquimbo"
$ws.Range("E114").Value2 = 4
$ws.Range("F114").Value2 = 100112044
$ws.Range("G114").Value2 = "Perejil"
$ws.Range("H114").Value2 = "Sin especificar"
$ws.Range("I114").Value2 = "Primera"
$ws.Range("N114").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O114").Value2 = "Provincia del Elquí"
$ws.Range("Q114").Value2 = 1.5
$ws.Range("R114").Value2 = "Hortaliza"
$ws.Range("D114").Value2 = 44306
$ws.Range("D114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J114").Value2 = 3400
$ws.Range("K114").Value2 = 2000
$ws.Range("L114").Value2 = 2500
$ws.Range("M114").Value2 = 2250
$ws.Range("P114").Value2 = 1500

$ws.Range("A115").Value2 = 8
$ws.Range("B115").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C115").Value2 = "Coquimbo"
$ws.Range("E115").Value2 = 4
$ws.Range("F115").Value2 = 100112044
$ws.Range("G115").Value2 = "Perejil"
$ws.Range("H115").Value2 = "Sin especificar"
$ws.Range("I115").Value2 = "Primera"
$ws.Range("N115").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O115").Value2 = "Provincia del Elquí"
$ws.Range("Q115").Value2 = 1.5
$ws.Range("R115").Value2 = "Hortaliza"
$ws.Range("D115").Value2 = 44189
$ws.Range("D115").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J115").Value2 = 2600
$ws.Range("K115").Value2 = 1400
$ws.Range("L115").Value2 = 1500
$ws.Range("M115").Value2 = 1450
$ws.Range("P115").Value2 = 967

$ws.Range("A116").Value2 = 8
$ws.Range("B116").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C116").Value2 = "Coquimbo"
$ws.Range("E116").Value2 = 4
$ws.Range("F116").Value2 = 100112044
$ws.Range("G116").Value2 = "Perejil"
$ws.Range("H116").Value2 = "Sin especificar"
$ws.Range("I116").Value2 = "Primera"
$ws.Range("N116").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O116").Value2 = "Provincia del Elquí"
$ws.Range("Q116").Value2 = 1.5
$ws.Range("R116").Value2 = "Hortaliza"
$ws.Range("D116").Value2 = 44299
$ws.Range("D116").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J116").Value2 = 3400
$ws.Range("K116").Value2 = 2000
$ws.Range("L116").Value2 = 2500
$ws.Range("M116").Value2 = 2250
$ws.Range("P116").Value2 = 1500

$ws.Range("A117").Value2 = 8
$ws.Range("B117").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C117").Value2 = "Coquimbo"
$ws.Range("E117").Value2 = 4
$ws.Range("F117").Value2 = 100112044
$ws.Range("G117").Value2 = "Perejil"
$ws.Range("H117").Value2 = "Sin especificar"
$ws.Range("I117").Value2 = "Primera"
$ws.Range("N117").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O117").Value2 = "Provincia del Elquí"
$ws.Range("Q117").Value2 = 1.5
$ws.Range("R117").Value2 = "Hortaliza"
$ws.Range("D117").Value2 = 44392
$ws.Range("D117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J117").Value2 = 3320
$ws.Range("K117").Value2 = 1500
$ws.Range("L117").Value2 = 2000
$ws.Range("M117").Value2 = 1750
$ws.Range("P117").Value2 = 1167

